$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.088.93'
$ws.Range('E2').Value = '  +2.38%  '
$ws.Range('D3').Value = '1.806.70'
$ws.Range('E3').Value = '  +0.39%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.007'
$ws.Range('E4').Value = '  +0.33%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '339.39'
$ws.Range('E5').Value = '  +0.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.004'
$ws.Range('E6').Value = '  +0.40%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3903'
$ws.Range('E7').Value = '  +2.49%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3473'
$ws.Range('E8').Value = '  +0.37%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '48.26'
$ws.Range('E9').Value = '  -0.99%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.189'
$ws.Range('E10').Value = '  -1.31%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07568'
$ws.Range('E11').Value = '  +0.47%  '
$ws.Range('E12').Value = '  +0.34%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.99'
$ws.Range('E13').Value = '  -0.68%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.492'
$ws.Range('E14').Value = '  -0.05%  '
$ws.Range('D15').Value = '1.813.96'
$ws.Range('E15').Value = '  +0.89%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.126'
$ws.Range('E16').Value = '  +0.34%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001099'
$ws.Range('E17').Value = '  -0.81%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06714'
$ws.Range('E18').Value = '  +0.85%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '84.99'
$ws.Range('E19').Value = '  -0.04%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.004'
$ws.Range('E20').Value = '  +0.39%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.72'
$ws.Range('E21').Value = '  +1.72%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.556'
$ws.Range('E22').Value = '  +0.23%  '
$ws.Range('D23').Value = '28.090.88'
$ws.Range('E23').Value = '  +2.42%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.40'
$ws.Range('E24').Value = '  -1.63%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.418'
$ws.Range('E25').Value = '  -0.25%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.496'
$ws.Range('E26').Value = '  -0.17%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.520'
$ws.Range('E27').Value = '  -2.06%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '21.23'
$ws.Range('E28').Value = '  -1.41%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '154.30'
$ws.Range('E29').Value = '  +1.24%  '
$ws.Range('D30').Value = '2.020.67'
$ws.Range('E30').Value = '  +0.95%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '135.58'
$ws.Range('E31').Value = '  +1.21%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.046'
$ws.Range('E32').Value = '  -0.09%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.135'
$ws.Range('E33').Value = '  -0.16%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08796'
$ws.Range('E34').Value = '  +1.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '12.96'
$ws.Range('E35').Value = '  -2.86%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06550'
$ws.Range('E36').Value = '  +2.30%  '
$ws.Range('B37').Value = 'InternetComputer(DFINITY)'
$ws.Range('C37').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.462'
$ws.Range('E37').Value = '  -0.12%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02423'
$ws.Range('E38').Value = '  +3.35%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.6907'
$ws.Range('E39').Value = '  -0.23%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.613'
$ws.Range('E40').Value = '  -1.99%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2214'
$ws.Range('E41').Value = '  +0.17%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.262'
$ws.Range('E42').Value = '  -1.13%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.505'
$ws.Range('E43').Value = '  -4.65%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.68'
$ws.Range('E44').Value = '  +1.55%  '
$ws.Range('B45').Value = 'Frax'
$ws.Range('C45').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.003'
$ws.Range('E45').Value = '  +0.34%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6422'
$ws.Range('E46').Value = '  -0.75%  '
$ws.Range('B47').Value = 'PancakeSwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.873'
$ws.Range('E47').Value = '  +0.11%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.151'
$ws.Range('E48').Value = '  +0.21%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '131.02'
$ws.Range('E49').Value = '  +0.31%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07188'
$ws.Range('E50').Value = '  -0.23%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '79.85'
$ws.Range('E51').Value = '  -0.31%  '
